$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5160.4
$ws.Range("I2").Value = 7291.4287
$ws.Range("K2").Value = 7291.4287
$ws.Range("M2").Value = -7178.4287
$ws.Range("H6").Value = 6668.3335
$ws.Range("I6").Value = 5
$ws.Range("K6").Value = 15
$ws.Range("M6").Value = 97
$ws.Range("H28").Value = 63348
$ws.Range("I28").Value = 73137.5
$ws.Range("J28").Value = 17663.666
$ws.Range("K28").Value = 73137.5
$ws.Range("L28").Value = 17663.666
$ws.Range("M28").Value = -72652.5
$ws.Range("N28").Value = -18633.666
$ws.Range("H33").Value = 15672183
$ws.Range("I33").Value = 57897.77
$ws.Range("J33").Value = 83334080
$ws.Range("K33").Value = 57897.77
$ws.Range("L33").Value = 83334080
$ws.Range("M33").Value = -57668.77
$ws.Range("N33").Value = -83334538
$ws.Range("H38").Value = 105.63636
$ws.Range("I38").Value = 105.63636
$ws.Range("K38").Value = 316.90908
$ws.Range("M38").Value = 55.09091999999998
$ws.Range("H40").Value = 3068.5789
$ws.Range("I40").Value = 2857.2856
$ws.Range("J40").Value = 3191.8333
$ws.Range("K40").Value = 2857.2856
$ws.Range("L40").Value = 3191.8333
$ws.Range("M40").Value = -2682.2856
$ws.Range("N40").Value = -3541.8333
$ws.Range("H41").Value = 1437233.2
$ws.Range("I41").Value = 5000056
$ws.Range("J41").Value = 12104.2
$ws.Range("K41").Value = 5000056
$ws.Range("L41").Value = 12104.2
$ws.Range("M41").Value = -4999616
$ws.Range("N41").Value = -12984.2
$ws.Range("H62").Value = 1205222.6
$ws.Range("I62").Value = 5150002
$ws.Range("J62").Value = 78142.86
$ws.Range("K62").Value = 5150002
$ws.Range("L62").Value = 78142.86
$ws.Range("M62").Value = -5149378
$ws.Range("N62").Value = -79390.86
$ws.Range("H65").Value = 1205222.6
$ws.Range("I65").Value = 5150002
$ws.Range("J65").Value = 78142.86
$ws.Range("K65").Value = 25750010
$ws.Range("L65").Value = 390714.3
$ws.Range("M65").Value = -25746890
$ws.Range("N65").Value = -396954.3
$ws.Range("H70").Value = 466005.1
$ws.Range("I70").Value = 1135192.5
$ws.Range("J70").Value = 2721.4614
$ws.Range("K70").Value = 3405577.5
$ws.Range("L70").Value = 8164.3842
$ws.Range("M70").Value = -3405307.5
$ws.Range("N70").Value = -8704.3842
$ws.Range("H73").Value = 466005.1
$ws.Range("I73").Value = 1135192.5
$ws.Range("J73").Value = 2721.4614
$ws.Range("K73").Value = 3405577.5
$ws.Range("L73").Value = 8164.3842
$ws.Range("M73").Value = -3404641.5
$ws.Range("N73").Value = -10036.3842
$ws.Range("H76").Value = 3997448.2
$ws.Range("I76").Value = 4107.4
$ws.Range("K76").Value = 4107.4
$ws.Range("M76").Value = -3792.4
$ws.Range("H79").Value = 3997448.2
$ws.Range("I79").Value = 4107.4
$ws.Range("K79").Value = 4107.4
$ws.Range("M79").Value = -3015.4
$ws.Range("H80").Value = 22727272
$ws.Range("I80").Value = 22727272
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 68181816
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -68180818
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 22727272
$ws.Range("I83").Value = 22727272
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 204545448
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -204540456
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 14327426
$ws.Range("I86").Value = 8993.5
$ws.Range("J86").Value = 20054800
$ws.Range("K86").Value = 8993.5
$ws.Range("L86").Value = 20054800
$ws.Range("M86").Value = -7870.5
$ws.Range("N86").Value = -20057046
$ws.Range("H89").Value = 14327426
$ws.Range("I89").Value = 8993.5
$ws.Range("J89").Value = 20054800
$ws.Range("K89").Value = 44967.5
$ws.Range("L89").Value = 100274000
$ws.Range("M89").Value = -39351.5
$ws.Range("N89").Value = -100285232
$ws.Range("H98").Value = 526.6667
$ws.Range("I98").Value = 526.6667
$ws.Range("K98").Value = 526.6667
$ws.Range("M98").Value = 971.3333
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 587.05
$ws.Range("I107").Value = 602.25
$ws.Range("J107").Value = 526.25
$ws.Range("K107").Value = 602.25
$ws.Range("L107").Value = 526.25
$ws.Range("M107").Value = 1317.75
$ws.Range("N107").Value = -4366.25
$ws.Range("H111").Value = 47666
$ws.Range("I111").Value = 31533.334
$ws.Range("J111").Value = 59765.5
$ws.Range("K111").Value = 94600.00199999999
$ws.Range("L111").Value = 179296.5
$ws.Range("M111").Value = -91533.00199999999
$ws.Range("N111").Value = -185430.5
$ws.Range("H113").Value = 71432630
$ws.Range("I113").Value = 111114320
$ws.Range("J113").Value = 5601.2
$ws.Range("K113").Value = 111114320
$ws.Range("L113").Value = 5601.2
$ws.Range("M113").Value = -111111066
$ws.Range("N113").Value = -12109.2
$ws.Range("H116").Value = 32469360
$ws.Range("I116").Value = 22823468
$ws.Range("J116").Value = 47627190
$ws.Range("K116").Value = 22823468
$ws.Range("L116").Value = 47627190
$ws.Range("M116").Value = -22820026
$ws.Range("N116").Value = -47634074
$ws.Range("H121").Value = 200363.73
$ws.Range("J121").Value = 203632.02
$ws.Range("L121").Value = 610896.0599999999
$ws.Range("N121").Value = -614390.0599999999
$ws.Range("H122").Value = 526.6667
$ws.Range("I122").Value = 526.6667
$ws.Range("K122").Value = 1580.0001
$ws.Range("M122").Value = 869.9999
$ws.Range("H125").Value = 909
$ws.Range("I125").Value = 533.3333
$ws.Range("J125").Value = 2036
$ws.Range("K125").Value = 4799.9997
$ws.Range("L125").Value = 18324
$ws.Range("M125").Value = -2339.9997
$ws.Range("N125").Value = -23244
$ws.Range("H129").Value = 1997
$ws.Range("J129").Value = 1997
$ws.Range("L129").Value = 5991
$ws.Range("N129").Value = -15991
$ws.Range("H131").Value = 11478.105
$ws.Range("I131").Value = 4864
$ws.Range("J131").Value = 14530.77
$ws.Range("K131").Value = 14592
$ws.Range("L131").Value = 43592.31
$ws.Range("M131").Value = -9552
$ws.Range("N131").Value = -53672.31
$ws.Range("H132").Value = 3442.7534
$ws.Range("I132").Value = 3256.0186
$ws.Range("K132").Value = 9768.0558
$ws.Range("M132").Value = -7238.0558
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H135").Value = 71429290
$ws.Range("I135").Value = 71429290
$ws.Range("K135").Value = 642863610
$ws.Range("M135").Value = -642861075
$ws.Range("H137").Value = 3704.76
$ws.Range("I137").Value = 2143.25
$ws.Range("J137").Value = 3840.5435
$ws.Range("K137").Value = 6429.75
$ws.Range("L137").Value = 11521.6305
$ws.Range("M137").Value = -3879.75
$ws.Range("N137").Value = -16621.6305
$ws.Range("H138").Value = 4601.193
$ws.Range("I138").Value = 2464
$ws.Range("J138").Value = 5587.59
$ws.Range("K138").Value = 7392
$ws.Range("L138").Value = 16762.77
$ws.Range("M138").Value = -2252
$ws.Range("N138").Value = -27042.77
$ws.Range("H141").Value = 1449.375
$ws.Range("I141").Value = 1449.375
$ws.Range("K141").Value = 4348.125
$ws.Range("M141").Value = 831.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1424.8572
$ws.Range("I2").Value = 1329
$ws.Range("K2").Value = 1329
$ws.Range("M2").Value = -1216
$ws.Range("H32").Value = 215162.31
$ws.Range("I32").Value = 234970.66
$ws.Range("J32").Value = 2222.75
$ws.Range("K32").Value = 234970.66
$ws.Range("L32").Value = 2222.75
$ws.Range("M32").Value = -234683.66
$ws.Range("N32").Value = -2796.75
$ws.Range("H45").Value = 3525.4285
$ws.Range("I45").Value = 3335.8
$ws.Range("J45").Value = 3999.5
$ws.Range("K45").Value = 3335.8
$ws.Range("L45").Value = 3999.5
$ws.Range("M45").Value = -2958.8
$ws.Range("N45").Value = -4753.5
$ws.Range("H61").Value = 3340.4666
$ws.Range("I61").Value = 2371.3635
$ws.Range("J61").Value = 6005.5
$ws.Range("K61").Value = 2371.3635
$ws.Range("L61").Value = 6005.5
$ws.Range("M61").Value = -2159.3635
$ws.Range("N61").Value = -6429.5
$ws.Range("H74").Value = 7339.5293
$ws.Range("I74").Value = 6525.077
$ws.Range("J74").Value = 9986.5
$ws.Range("K74").Value = 6525.077
$ws.Range("L74").Value = 9986.5
$ws.Range("M74").Value = -5651.077
$ws.Range("N74").Value = -11734.5
$ws.Range("H77").Value = 7339.5293
$ws.Range("I77").Value = 6525.077
$ws.Range("J77").Value = 9986.5
$ws.Range("K77").Value = 32625.385
$ws.Range("L77").Value = 49932.5
$ws.Range("M77").Value = -28257.385
$ws.Range("N77").Value = -58668.5
$ws.Range("H97").Value = 355.70587
$ws.Range("I97").Value = 375.4375
$ws.Range("K97").Value = 375.4375
$ws.Range("M97").Value = 120.5625
$ws.Range("H110").Value = 90924220
$ws.Range("I110").Value = 111112720
$ws.Range("J110").Value = 76000
$ws.Range("K110").Value = 111112720
$ws.Range("L110").Value = 76000
$ws.Range("M110").Value = -111110675
$ws.Range("N110").Value = -80090
$ws.Range("H116").Value = 1424.8572
$ws.Range("I116").Value = 1329
$ws.Range("K116").Value = 1329
$ws.Range("M116").Value = 965
$ws.Range("H122").Value = 13891475
$ws.Range("I122").Value = 18520494
$ws.Range("K122").Value = 55561482
$ws.Range("M122").Value = -55559032
$ws.Range("H132").Value = 32261600
$ws.Range("I132").Value = 38465036
$ws.Range("J132").Value = 3725.2
$ws.Range("K132").Value = 115395108
$ws.Range("L132").Value = 11175.6
$ws.Range("M132").Value = -115392578
$ws.Range("N132").Value = -16235.6
$ws.Range("H136").Value = 3340.4666
$ws.Range("I136").Value = 2371.3635
$ws.Range("J136").Value = 6005.5
$ws.Range("K136").Value = 7114.0905
$ws.Range("L136").Value = 18016.5
$ws.Range("M136").Value = -4564.0905
$ws.Range("N136").Value = -23116.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1424.8572
$ws.Range("I3").Value = 1329
$ws.Range("K3").Value = 1329
$ws.Range("M3").Value = -1215
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -217
$ws.Range("H35").Value = 62055.25
$ws.Range("J35").Value = 62055.25
$ws.Range("L35").Value = 62055.25
$ws.Range("N35").Value = -62675.25
$ws.Range("H75").Value = 10242.4
$ws.Range("I75").Value = 10242.4
$ws.Range("K75").Value = 10242.4
$ws.Range("M75").Value = -9306.4
$ws.Range("H78").Value = 10242.4
$ws.Range("I78").Value = 10242.4
$ws.Range("K78").Value = 30727.2
$ws.Range("M78").Value = -26047.2
$ws.Range("H82").Value = 40828.332
$ws.Range("J82").Value = 52187
$ws.Range("L82").Value = 52187
$ws.Range("N82").Value = -52953
$ws.Range("H85").Value = 40828.332
$ws.Range("J85").Value = 52187
$ws.Range("L85").Value = 52187
$ws.Range("N85").Value = -54839
$ws.Range("H86").Value = 35716910
$ws.Range("I86").Value = 83335480
$ws.Range("J86").Value = 2980.875
$ws.Range("K86").Value = 83335480
$ws.Range("L86").Value = 2980.875
$ws.Range("M86").Value = -83334357
$ws.Range("N86").Value = -5226.875
$ws.Range("H89").Value = 35716910
$ws.Range("I89").Value = 83335480
$ws.Range("J89").Value = 2980.875
$ws.Range("K89").Value = 416677400
$ws.Range("L89").Value = 14904.375
$ws.Range("M89").Value = -416671784
$ws.Range("N89").Value = -26136.375
$ws.Range("H94").Value = 35719428
$ws.Range("I94").Value = 50004000
$ws.Range("K94").Value = 50004000
$ws.Range("M94").Value = -50003549
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 2301
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 2602
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 2602
$ws.Range("M105").Value = -253
$ws.Range("N105").Value = -6096
$ws.Range("H107").Value = 27810878
$ws.Range("I107").Value = 16545.908
$ws.Range("K107").Value = 16545.908
$ws.Range("M107").Value = -14625.908
$ws.Range("H117").Value = 70000
$ws.Range("J117").Value = 70000
$ws.Range("L117").Value = 70000
$ws.Range("N117").Value = -79178
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
$ws.Range("H134").Value = 2822.2222
$ws.Range("I134").Value = 2471.0417
$ws.Range("J134").Value = 5631.6665
$ws.Range("K134").Value = 7413.125100000001
$ws.Range("L134").Value = 16894.9995
$ws.Range("M134").Value = -4878.125100000001
$ws.Range("N134").Value = -21964.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2247.0715
$ws.Range("I22").Value = 373.25
$ws.Range("J22").Value = 2996.6
$ws.Range("K22").Value = 373.25
$ws.Range("L22").Value = 2996.6
$ws.Range("M22").Value = -23.25
$ws.Range("N22").Value = -3696.6
$ws.Range("H31").Value = 3827.842
$ws.Range("I31").Value = 1109.2
$ws.Range("J31").Value = 4239.758
$ws.Range("K31").Value = 1109.2
$ws.Range("L31").Value = 4239.758
$ws.Range("M31").Value = -814.2
$ws.Range("N31").Value = -4829.758
$ws.Range("H34").Value = 3827.842
$ws.Range("I34").Value = 1109.2
$ws.Range("J34").Value = 4239.758
$ws.Range("K34").Value = 1109.2
$ws.Range("L34").Value = 4239.758
$ws.Range("M34").Value = -907.2
$ws.Range("N34").Value = -4643.758
$ws.Range("H35").Value = 25005076
$ws.Range("I35").Value = 1120
$ws.Range("J35").Value = 66678332
$ws.Range("K35").Value = 1120
$ws.Range("L35").Value = 66678332
$ws.Range("M35").Value = -826
$ws.Range("N35").Value = -66678920
$ws.Range("H58").Value = 503207.3
$ws.Range("I58").Value = 1616.7142
$ws.Range("K58").Value = 1616.7142
$ws.Range("M58").Value = -1413.7142
$ws.Range("H62").Value = 15799.8
$ws.Range("H64").Value = 30270.5
$ws.Range("J64").Value = 30270.5
$ws.Range("L64").Value = 30270.5
$ws.Range("N64").Value = -30766.5
$ws.Range("H65").Value = 15799.8
$ws.Range("H67").Value = 30270.5
$ws.Range("J67").Value = 30270.5
$ws.Range("L67").Value = 30270.5
$ws.Range("N67").Value = -31986.5
$ws.Range("H94").Value = 1663.6666
$ws.Range("J94").Value = 1871.25
$ws.Range("L94").Value = 1871.25
$ws.Range("N94").Value = -2773.25
$ws.Range("H99").Value = 2383.7144
$ws.Range("I99").Value = 2118.0588
$ws.Range("J99").Value = 2794.2727
$ws.Range("K99").Value = 2118.0588
$ws.Range("L99").Value = 2794.2727
$ws.Range("M99").Value = -620.0587999999998
$ws.Range("N99").Value = -5790.2727
$ws.Range("H105").Value = 1869.75
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1130.4375
$ws.Range("I107").Value = 1130.4375
$ws.Range("K107").Value = 1130.4375
$ws.Range("M107").Value = 789.5625
$ws.Range("H111").Value = 80000
$ws.Range("J111").Value = 80000
$ws.Range("L111").Value = 80000
$ws.Range("N111").Value = -88180
$ws.Range("H122").Value = 3598.6
$ws.Range("I122").Value = 3598.6
$ws.Range("K122").Value = 10795.8
$ws.Range("M122").Value = -8345.799999999999
$ws.Range("H126").Value = 2383.7144
$ws.Range("I126").Value = 2118.0588
$ws.Range("J126").Value = 2794.2727
$ws.Range("K126").Value = 6354.176399999999
$ws.Range("L126").Value = 8382.8181
$ws.Range("M126").Value = -3884.176399999999
$ws.Range("N126").Value = -13322.8181
$ws.Range("H132").Value = 1431281.9
$ws.Range("I132").Value = 1002294.7
$ws.Range("J132").Value = 2503749.8
$ws.Range("K132").Value = 3006884.1
$ws.Range("L132").Value = 7511249.399999999
$ws.Range("M132").Value = -3004354.1
$ws.Range("N132").Value = -7516309.399999999
$ws.Range("H134").Value = 3462.2
$ws.Range("I134").Value = 2968.3333
$ws.Range("J134").Value = 4203
$ws.Range("K134").Value = 8904.999899999999
$ws.Range("L134").Value = 12609
$ws.Range("M134").Value = -6369.999899999999
$ws.Range("N134").Value = -17679
$ws.Range("H136").Value = 503207.3
$ws.Range("I136").Value = 1616.7142
$ws.Range("K136").Value = 4850.142599999999
$ws.Range("M136").Value = -2300.142599999999
$ws.Range("H141").Value = 101584.5
$ws.Range("I141").Value = 80000
$ws.Range("J141").Value = 103244.84
$ws.Range("K141").Value = 80000
$ws.Range("L141").Value = 103244.84
$ws.Range("M141").Value = -74820
$ws.Range("N141").Value = -113604.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20491330
$ws.Range("I4").Value = 21131526
$ws.Range("J4").Value = 5002
$ws.Range("K4").Value = 63394578
$ws.Range("L4").Value = 15006
$ws.Range("M4").Value = -63394466
$ws.Range("N4").Value = -15230
$ws.Range("H5").Value = 644.86957
$ws.Range("J5").Value = 836.5833
$ws.Range("L5").Value = 2509.7499
$ws.Range("N5").Value = -2733.7499
$ws.Range("H9").Value = 10248.4
$ws.Range("I9").Value = 14497.286
$ws.Range("J9").Value = 334.33334
$ws.Range("K9").Value = 43491.858
$ws.Range("L9").Value = 1003.00002
$ws.Range("M9").Value = -43267.858
$ws.Range("N9").Value = -1451.00002
$ws.Range("H10").Value = 712.4
$ws.Range("I10").Value = 310.2857
$ws.Range("K10").Value = 930.8571000000001
$ws.Range("M10").Value = -791.8571000000001
$ws.Range("H11").Value = 163694
$ws.Range("I11").Value = 52176.1
$ws.Range("K11").Value = 156528.3
$ws.Range("M11").Value = -156388.3
$ws.Range("H26").Value = 44
$ws.Range("I26").Value = 44
$ws.Range("K26").Value = 132
$ws.Range("M26").Value = 156
$ws.Range("H34").Value = 2034.625
$ws.Range("J34").Value = 3399.4443
$ws.Range("L34").Value = 10198.3329
$ws.Range("N34").Value = -10366.3329
$ws.Range("H39").Value = 3741.7896
$ws.Range("J39").Value = 3894.1667
$ws.Range("L39").Value = 11682.5001
$ws.Range("N39").Value = -12270.5001
$ws.Range("H40").Value = 76.77778000000001
$ws.Range("I40").Value = 81.92308
$ws.Range("J40").Value = 63.4
$ws.Range("K40").Value = 327.69232
$ws.Range("L40").Value = 253.6
$ws.Range("M40").Value = -258.69232
$ws.Range("N40").Value = -391.6
$ws.Range("H55").Value = 1672
$ws.Range("J55").Value = 2438.8
$ws.Range("L55").Value = 7316.400000000001
$ws.Range("N55").Value = -7670.400000000001
$ws.Range("H68").Value = 1316.25
$ws.Range("I68").Value = 1117.7273
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 3353.1819
$ws.Range("L68").Value = 10500
$ws.Range("M68").Value = -2542.1819
$ws.Range("N68").Value = -12122
$ws.Range("H71").Value = 1316.25
$ws.Range("I71").Value = 1117.7273
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 10059.5457
$ws.Range("L71").Value = 31500
$ws.Range("M71").Value = -6003.545700000001
$ws.Range("N71").Value = -39612
$ws.Range("H107").Value = 1850.409
$ws.Range("I107").Value = 1615.6428
$ws.Range("J107").Value = 2261.25
$ws.Range("K107").Value = 4846.928400000001
$ws.Range("L107").Value = 6783.75
$ws.Range("M107").Value = -2926.928400000001
$ws.Range("N107").Value = -10623.75
$ws.Range("H113").Value = 656.3871
$ws.Range("I113").Value = 240.85715
$ws.Range("J113").Value = 998.58826
$ws.Range("K113").Value = 722.5714499999999
$ws.Range("L113").Value = 2995.76478
$ws.Range("M113").Value = 1447.42855
$ws.Range("N113").Value = -7335.76478
$ws.Range("H119").Value = 11347.261
$ws.Range("I119").Value = 2989
$ws.Range("K119").Value = 8967
$ws.Range("M119").Value = -4129
$ws.Range("H121").Value = 123515.89
$ws.Range("I121").Value = 883
$ws.Range("J121").Value = 184832.33
$ws.Range("K121").Value = 2649
$ws.Range("L121").Value = 554496.99
$ws.Range("M121").Value = -1339
$ws.Range("N121").Value = -557116.99
$ws.Range("H122").Value = 358.18182
$ws.Range("J122").Value = 381.7143
$ws.Range("L122").Value = 3435.4287
$ws.Range("N122").Value = -8335.4287
$ws.Range("H128").Value = 214666
$ws.Range("I128").Value = 214666
$ws.Range("K128").Value = 643998
$ws.Range("M128").Value = -639018
$ws.Range("H129").Value = 1993.0714
$ws.Range("I129").Value = 1377
$ws.Range("K129").Value = 4131
$ws.Range("M129").Value = 869
$ws.Range("H131").Value = 15271.75
$ws.Range("J131").Value = 15271.75
$ws.Range("L131").Value = 45815.25
$ws.Range("N131").Value = -55895.25
$ws.Range("H132").Value = 1493
$ws.Range("I132").Value = 1470
$ws.Range("K132").Value = 13230
$ws.Range("M132").Value = -10700
$ws.Range("H133").Value = 2765
$ws.Range("I133").Value = 2765
$ws.Range("K133").Value = 8295
$ws.Range("M133").Value = -3235
$ws.Range("H134").Value = 2377.875
$ws.Range("I134").Value = 2377.875
$ws.Range("K134").Value = 7133.625
$ws.Range("M134").Value = -2063.625
$ws.Range("H135").Value = 644.86957
$ws.Range("J135").Value = 836.5833
$ws.Range("L135").Value = 7529.2497
$ws.Range("N135").Value = -12599.2497
$ws.Range("H136").Value = 4178.75
$ws.Range("I136").Value = 4905
$ws.Range("K136").Value = 14715
$ws.Range("M136").Value = -9615
$ws.Range("H137").Value = 2614.2778
$ws.Range("J137").Value = 2939.9333
$ws.Range("L137").Value = 8819.7999
$ws.Range("N137").Value = -19019.7999
$ws.Range("H140").Value = 1093.0769
$ws.Range("I140").Value = 1093.0769
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3279.2307
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 1900.7693
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1616.7693
$ws.Range("I2").Value = 1324.125
$ws.Range("J2").Value = 2085
$ws.Range("K2").Value = 1324.125
$ws.Range("L2").Value = 2085
$ws.Range("M2").Value = -1211.125
$ws.Range("N2").Value = -2311
$ws.Range("H51").Value = 80000
$ws.Range("J51").Value = 80000
$ws.Range("L51").Value = 80000
$ws.Range("N51").Value = -81018
$ws.Range("H70").Value = 13114.875
$ws.Range("I70").Value = 14978
$ws.Range("J70").Value = 11251.75
$ws.Range("K70").Value = 14978
$ws.Range("L70").Value = 11251.75
$ws.Range("M70").Value = -14708
$ws.Range("N70").Value = -11791.75
$ws.Range("H73").Value = 13114.875
$ws.Range("I73").Value = 14978
$ws.Range("J73").Value = 11251.75
$ws.Range("K73").Value = 14978
$ws.Range("L73").Value = 11251.75
$ws.Range("M73").Value = -14042
$ws.Range("N73").Value = -13123.75
$ws.Range("H80").Value = 3594.3635
$ws.Range("I80").Value = 3539.8572
$ws.Range("J80").Value = 3689.75
$ws.Range("K80").Value = 3539.8572
$ws.Range("L80").Value = 3689.75
$ws.Range("M80").Value = -2541.8572
$ws.Range("N80").Value = -5685.75
$ws.Range("H83").Value = 3594.3635
$ws.Range("I83").Value = 3539.8572
$ws.Range("J83").Value = 3689.75
$ws.Range("K83").Value = 17699.286
$ws.Range("L83").Value = 18448.75
$ws.Range("M83").Value = -12707.286
$ws.Range("N83").Value = -28432.75
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492
$ws.Range("H97").Value = 1914.9131
$ws.Range("I97").Value = 1964.6111
$ws.Range("K97").Value = 1964.6111
$ws.Range("M97").Value = -1468.6111
$ws.Range("H99").Value = 9950
$ws.Range("I99").Value = 9950
$ws.Range("K99").Value = 9950
$ws.Range("M99").Value = -7704
$ws.Range("H102").Value = 2561.8667
$ws.Range("J102").Value = 3999.75
$ws.Range("L102").Value = 3999.75
$ws.Range("N102").Value = -7243.75
$ws.Range("H113").Value = 2682
$ws.Range("I113").Value = 1872.4667
$ws.Range("J113").Value = 3785.9092
$ws.Range("K113").Value = 1872.4667
$ws.Range("L113").Value = 3785.9092
$ws.Range("M113").Value = 297.5333000000001
$ws.Range("N113").Value = -8125.9092
$ws.Range("H122").Value = 100002500
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 250002500
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 750007500
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -750012400
$ws.Range("H132").Value = 175508.78
$ws.Range("I132").Value = 240922.47
$ws.Range("J132").Value = 3797.8125
$ws.Range("K132").Value = 722767.41
$ws.Range("L132").Value = 11393.4375
$ws.Range("M132").Value = -720237.41
$ws.Range("N132").Value = -16453.4375
$ws.Range("H135").Value = 145159.8
$ws.Range("I135").Value = 108000
$ws.Range("J135").Value = 154449.75
$ws.Range("K135").Value = 108000
$ws.Range("L135").Value = 154449.75
$ws.Range("M135").Value = -102930
$ws.Range("N135").Value = -164589.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31253718
$ws.Range("I7").Value = 71431430
$ws.Range("J7").Value = 4388.778
$ws.Range("K7").Value = 71431430
$ws.Range("L7").Value = 4388.778
$ws.Range("M7").Value = -71431318
$ws.Range("N7").Value = -4612.778
$ws.Range("H22").Value = 3587929.5
$ws.Range("I22").Value = 4450
$ws.Range("J22").Value = 6454713
$ws.Range("K22").Value = 4450
$ws.Range("L22").Value = 6454713
$ws.Range("M22").Value = -4155
$ws.Range("N22").Value = -6455303
$ws.Range("H27").Value = 3587929.5
$ws.Range("I27").Value = 4450
$ws.Range("J27").Value = 6454713
$ws.Range("K27").Value = 4450
$ws.Range("L27").Value = 6454713
$ws.Range("M27").Value = -4343
$ws.Range("N27").Value = -6454927
$ws.Range("H40").Value = 3257.7144
$ws.Range("I40").Value = 3257.7144
$ws.Range("K40").Value = 3257.7144
$ws.Range("M40").Value = -3121.7144
$ws.Range("H46").Value = 2945.3333
$ws.Range("I46").Value = 2416.6667
$ws.Range("J46").Value = 3156.8
$ws.Range("K46").Value = 2416.6667
$ws.Range("L46").Value = 3156.8
$ws.Range("M46").Value = -2228.6667
$ws.Range("N46").Value = -3532.8
$ws.Range("H55").Value = 415.58334
$ws.Range("I55").Value = 400.6316
$ws.Range("J55").Value = 472.4
$ws.Range("K55").Value = 400.6316
$ws.Range("L55").Value = 472.4
$ws.Range("M55").Value = -227.6316
$ws.Range("N55").Value = -818.4
$ws.Range("H58").Value = 200001650
$ws.Range("I58").Value = 2055.5
$ws.Range("K58").Value = 2055.5
$ws.Range("M58").Value = -1795.5
$ws.Range("H61").Value = 233391.81
$ws.Range("I61").Value = 276494.3
$ws.Range("J61").Value = 5564.2856
$ws.Range("K61").Value = 276494.3
$ws.Range("L61").Value = 5564.2856
$ws.Range("M61").Value = -276292.3
$ws.Range("N61").Value = -5968.2856
$ws.Range("H68").Value = 29167.666
$ws.Range("I68").Value = 3750
$ws.Range("K68").Value = 3750
$ws.Range("M68").Value = -3001
$ws.Range("H71").Value = 29167.666
$ws.Range("I71").Value = 3750
$ws.Range("K71").Value = 18750
$ws.Range("M71").Value = -15006
$ws.Range("H82").Value = 4025.5715
$ws.Range("I82").Value = 4025.5715
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 4025.5715
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3664.5715
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 4025.5715
$ws.Range("I85").Value = 4025.5715
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 4025.5715
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2777.5715
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 7500
$ws.Range("I93").Value = 10000
$ws.Range("K93").Value = 10000
$ws.Range("M93").Value = -8752
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 233391.81
$ws.Range("I113").Value = 276494.3
$ws.Range("J113").Value = 5564.2856
$ws.Range("K113").Value = 276494.3
$ws.Range("L113").Value = 5564.2856
$ws.Range("M113").Value = -274324.3
$ws.Range("N113").Value = -9904.285599999999
$ws.Range("H122").Value = 3740
$ws.Range("I122").Value = 2918.182
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 8754.545999999998
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -6304.545999999998
$ws.Range("N122").Value = -22900
$ws.Range("H126").Value = 31253718
$ws.Range("I126").Value = 71431430
$ws.Range("J126").Value = 4388.778
$ws.Range("K126").Value = 214294290
$ws.Range("L126").Value = 13166.334
$ws.Range("M126").Value = -214291820
$ws.Range("N126").Value = -18106.334
$ws.Range("H132").Value = 6582
$ws.Range("I132").Value = 3544.818
$ws.Range("J132").Value = 9619.182000000001
$ws.Range("K132").Value = 10634.454
$ws.Range("L132").Value = 28857.546
$ws.Range("M132").Value = -8104.454000000002
$ws.Range("N132").Value = -33917.546
$ws.Range("H136").Value = 5514.9546
$ws.Range("I136").Value = 5053.1714
$ws.Range("J136").Value = 7310.778
$ws.Range("K136").Value = 15159.5142
$ws.Range("L136").Value = 21932.334
$ws.Range("M136").Value = -12609.5142
$ws.Range("N136").Value = -27032.334
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 66497.5
$ws.Range("J47").Value = 69000
$ws.Range("L47").Value = 69000
$ws.Range("N47").Value = -70144
$ws.Range("H62").Value = 51197450
$ws.Range("J62").Value = 55563024
$ws.Range("L62").Value = 55563024
$ws.Range("N62").Value = -55564272
$ws.Range("H65").Value = 51197450
$ws.Range("J65").Value = 55563024
$ws.Range("L65").Value = 277815120
$ws.Range("N65").Value = -277821360
$ws.Range("H81").Value = 15392259
$ws.Range("J81").Value = 18189952
$ws.Range("L81").Value = 36379904
$ws.Range("N81").Value = -36382026
$ws.Range("H82").Value = 13749.5
$ws.Range("I82").Value = 7500
$ws.Range("J82").Value = 19999
$ws.Range("K82").Value = 7500
$ws.Range("L82").Value = 19999
$ws.Range("M82").Value = -7117
$ws.Range("N82").Value = -20765
$ws.Range("H84").Value = 15392259
$ws.Range("J84").Value = 18189952
$ws.Range("L84").Value = 181899520
$ws.Range("N84").Value = -181910128
$ws.Range("H85").Value = 13749.5
$ws.Range("I85").Value = 7500
$ws.Range("J85").Value = 19999
$ws.Range("K85").Value = 7500
$ws.Range("L85").Value = 19999
$ws.Range("M85").Value = -6174
$ws.Range("N85").Value = -22651
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 9355.200000000001
$ws.Range("I96").Value = 8388.5
$ws.Range("J96").Value = 9999.666999999999
$ws.Range("K96").Value = 8388.5
$ws.Range("L96").Value = 9999.666999999999
$ws.Range("M96").Value = -7015.5
$ws.Range("N96").Value = -12745.667
$ws.Range("H100").Value = 722.7143
$ws.Range("I100").Value = 1063.5
$ws.Range("J100").Value = 268.33334
$ws.Range("K100").Value = 2127
$ws.Range("L100").Value = 536.66668
$ws.Range("M100").Value = -1586
$ws.Range("N100").Value = -1618.66668
$ws.Range("H107").Value = 759.8
$ws.Range("I107").Value = 499.66666
$ws.Range("K107").Value = 1498.99998
$ws.Range("M107").Value = 421.0000199999999
$ws.Range("H122").Value = 1798.9445
$ws.Range("I122").Value = 1687.4
$ws.Range("J122").Value = 2356.6667
$ws.Range("K122").Value = 5062.200000000001
$ws.Range("L122").Value = 7070.000100000001
$ws.Range("M122").Value = -2612.200000000001
$ws.Range("N122").Value = -11970.0001
$ws.Range("H132").Value = 388308.56
$ws.Range("I132").Value = 574579.0600000001
$ws.Range("J132").Value = 4810.5293
$ws.Range("K132").Value = 1723737.18
$ws.Range("L132").Value = 14431.5879
$ws.Range("M132").Value = -1721207.18
$ws.Range("N132").Value = -19491.5879
$ws.Range("H136").Value = 4043.4546
$ws.Range("I136").Value = 2905.5454
$ws.Range("J136").Value = 5181.364
$ws.Range("K136").Value = 8716.636200000001
$ws.Range("L136").Value = 15544.092
$ws.Range("M136").Value = -6166.636200000001
$ws.Range("N136").Value = -20644.092
$ws.Range("H140").Value = 46684.4
$ws.Range("J140").Value = 46684.4
$ws.Range("L140").Value = 46684.4
$ws.Range("N140").Value = -57044.4
$ws.Range("H141").Value = 42499.25
$ws.Range("J141").Value = 42499.25
$ws.Range("L141").Value = 42499.25
$ws.Range("N141").Value = -52859.25
